$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header column: Colesterol (placed after Creatinina, before Tiroxina Livre)
$ws.Range("Y1").Value = "Colesterol"

# Data entered in typed (pre-sort) order: Vanessa, Davi, Nathalia, Olivio
# row 2 -> Vanessa
$ws.Range("A2").Value = "VANESSA HINSELMANN DOS SANTOS"
$ws.Range("E2").Value = "35,90"
$ws.Range("D2").Value = "12,00"
$ws.Range("J2").Value = 7.46
$ws.Range("N2").Value = 27
$ws.Range("S2").Value = 10
$ws.Range("T2").Value = 238
$ws.Range("Y2").Value = 207

# row 3 -> Davi
$ws.Range("A3").Value = "DAVI KOBUS ZORZI"
$ws.Range("E3").Value = "39,4"
$ws.Range("D3").Value = "13,5"
$ws.Range("J3").Value = 100
$ws.Range("N3").Value = 52
$ws.Range("S3").Value = 4
$ws.Range("T3").Value = 474

# row 4 -> Nathalia
$ws.Range("A4").Value = "NATHALIA MACENA CUSTODIO"
$ws.Range("E4").Value = "36,6"
$ws.Range("D4").Value = "12,5"
$ws.Range("S4").Value = "8,4"
$ws.Range("N4").Value = "35,0"
$ws.Range("Y4").Value = "162,8"
$ws.Range("J4").Value = 3.69
$ws.Range("T4").Value = 243

# row 5 -> Olivio
$ws.Range("A5").Value = "OLIVIO JUVENCIO DA SILVA"
$ws.Range("E5").Value = "41,9"
$ws.Range("D5").Value = "14,7"
$ws.Range("J5").Value = 6.78
$ws.Range("N5").Value = 66
$ws.Range("S5").Value = 7
$ws.Range("Y5").Value = 260

# Sort the entered rows alphabetically by patient name (column A)
$sortRange = $ws.Range("A2:AC5")
$sortRange.Sort($ws.Range("A2:A5"), 1)

# Widen columns to fit the new content (Patient name column, and columns
# that were resized alongside it)
$ws.Columns("A").ColumnWidth = 33.33333333333333
$ws.Columns("L").ColumnWidth = 14.333333333333332
$ws.Columns("M").ColumnWidth = 9.666666666666666

# Final selection left by the user
$ws.Range("F6").Select()
